# Applies the cryptos-list price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.211.47'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.828.33'
$ws.Range('D4').Value = '''0.9992'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''236.68'
$ws.Range('E5').Value = '  -1.34%  '
$ws.Range('D6').Value = '''0.6081'
$ws.Range('E6').Value = '  -3.48%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D9').Value = '''0.2813'
$ws.Range('E9').Value = '  -2.73%  '
$ws.Range('E10').Value = '  -4.73%  '
$ws.Range('D11').Value = '''0.07665'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '1.823.05'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D13').Value = '''4.826'
$ws.Range('E13').Value = '  -2.91%  '
$ws.Range('D14').Value = '''0.00001010'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').Value = '''0.6346'
$ws.Range('E15').Value = '  -6.05%  '
$ws.Range('D16').Value = '2.071.87'
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('E17').Value = '  -3.11%  '
$ws.Range('D18').Value = '''5.889'
$ws.Range('E18').Value = '  -5.55%  '
$ws.Range('D19').Value = '29.191.95'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '''228.32'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  -4.13%  '
$ws.Range('D23').Value = '''7.029'
$ws.Range('E23').Value = '  -4.54%  '
$ws.Range('D24').Value = '''1.000'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = '''154.68'
$ws.Range('E25').Value = '  -2.09%  '
$ws.Range('D26').Value = '''8.076'
$ws.Range('E26').Value = '  -4.91%  '
$ws.Range('D27').Value = '''0.1304'
$ws.Range('E27').Value = '  -3.22%  '
$ws.Range('D28').Value = '''16.59'
$ws.Range('E28').Value = '  -4.86%  '
$ws.Range('D29').Value = '''1.483'
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('D30').Value = '''0.06472'
$ws.Range('E30').Value = '  -6.67%  '
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('D32').Value = '''3.827'
$ws.Range('E32').Value = '  -5.32%  '
$ws.Range('E33').Value = '  -5.98%  '
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('D35').Value = '''1.753'
$ws.Range('E35').Value = '  -3.94%  '
$ws.Range('D36').Value = '''0.6511'
$ws.Range('E36').Value = '  -6.80%  '
$ws.Range('D37').Value = '''2.545'
$ws.Range('E37').Value = '  -1.60%  '
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('D39').Value = '1.218.60'
$ws.Range('E39').Value = '  -1.23%  '
$ws.Range('D40').Value = '''0.01754'
$ws.Range('E40').Value = '  -4.86%  '
$ws.Range('D41').Value = '''6.554'
$ws.Range('E41').Value = '  -3.56%  '
$ws.Range('D42').Value = '''0.9326'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').Value = '''0.9995'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '''101.09'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = '1.982.90'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''63.14'
$ws.Range('E46').Value = '  -3.25%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '''0.00000000119'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''8.619'
$ws.Range('E48').Value = '  -3.24%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''1.614'
$ws.Range('E49').Value = '  -5.22%  '
$ws.Range('E50').Value = '  -5.47%  '
$ws.Range('E51').Value = '  -2.64%  '
